# EdisonLogs/gps.xlsx — "handler for pop up in chrome"
#
# Sheet1's place list is trimmed down to a single new entry ("Duren"),
# Sheet1 becomes the active/selected tab (it previously was Sheet2), and
# both sheets' selections move to reflect the new UI state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: drop the old Madrid/bilbao/valencia rows, add "Duren" ---
$ws1.Range("A2:C4").Clear()
$ws1.Range("A2").Value = "Duren"

# Column B on Sheet1 re-fits to its (now much shorter) content.
$ws1.Columns.Item(2).ColumnWidth = 20.71

# --- Selections / active tab ---
# Sheet2's selection moves off the old C22:C31 block to F22, and it stops
# being the active tab...
$ws2.Range("F22").Select()

# ...while Sheet1 becomes the active tab with C2 selected.
$ws1.Select()
$ws1.Range("C2").Select()
